$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 56: phone number was stored as text "09876543" (leading zero);
# convert it to the numeric value 9876543 to match the rest of the sheet. ---
$ws.Range("A56").Value = 9876543

# --- Append new payment row 57: 09876543 (Cash) 2025-08-18T18:05:43 ---
# Column A keeps the phone number's leading zero, so it must stay TEXT
# (not auto-coerced to the number 9876543 like a normal Value assignment
# would do). Force text via NumberFormat, assign, then restore the cell's
# style to Normal so no stray formatting is left behind on the cell.
$aCell = $ws.Range("A57")
$aCell.NumberFormat = "@"
$aCell.Value = "09876543"
$aCell.Style = "Normal"

$ws.Range("C57").Value = "Cash"
$ws.Range("D57").Value = "2025-08-18T18:05:43"
$ws.Range("E57").Value = 120
$ws.Range("G57").Value = 120
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
